# Update "Generate Report for Handback" timestamps on the per-language
# handback status sheets ("zh-cn" and "de-de").
#
# Column D = Correspond Handoff Datetime
# Column G = Correspond Handback DateTime
# Row 2 on each sheet holds the relevant record.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-17 09:45:49"
$wsZhCn.Range("G2").Value = "2016-02-17 09:46:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-17 09:46:00"
$wsDeDe.Range("G2").Value = "2016-02-17 09:47:12"
